# signal connections.xlsx - rework SPI / servo PWM pin assignments on the
# "XMOS Connections" sheet so SPI lands on the boot-from-flash pins.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS Connections")

# --- Column width tweaks: two new "in-between" columns picked up real widths ---
$ws.Columns.Item(7).ColumnWidth = 16.86    # column G ~17.71 stored width
$ws.Columns.Item(13).ColumnWidth = 13.86   # column M ~14.71 stored width

# --- Re-assigned signal names ---
$ws.Range("E2").Value2 = "SPI_MISO"
$ws.Range("E3").Value2 = "EEPROM_SS"

$ws.Range("N8").ClearContents()
$ws.Range("N9").ClearContents()

$ws.Range("E12").Value2 = "SPI_SCK"
$ws.Range("E13").Value2 = "SPI_MOSI"

$ws.Range("K14").Value2 = "SERVO_1"

$ws.Range("K24").Value2 = "SERVO_2"

$ws.Range("E26").Value2 = "GYRO_SS"
$ws.Range("K26").Value2 = "SERVO_0"

$ws.Range("E36").Value2 = "uSD_SS"
$ws.Range("K36").Value2 = "SERVO_3"

$ws.Range("E37").ClearContents()
$ws.Range("K37").Value2 = "SERVO_4"

$ws.Range("K38").Value2 = "SERVO_5"
$ws.Range("K39").Value2 = "PWM_0"
$ws.Range("K40").Value2 = "PWM_1"
$ws.Range("K41").Value2 = "PWM_2"

# --- Selection moved to the newly re-pinned servo block ---
$ws.Range("M36:M41").Select() | Out-Null
